$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three condition probability values
$ws.Range("C14").Value = 0.05
$ws.Range("C20").Value = 0.05
$ws.Range("C42").Value = 0.05

# Update the view: scroll position and active cell selection
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
